# Update the presentation per the target diff:
#  - refresh every auto-updating "date" placeholder (slide master, every
#    custom layout, and the notes master) from 4/28/21 -> 5/2/21
#  - tweak the copy on the two "Automate.../Safely promote..." callouts
#    on slide 1

$p = $ppt.ActivePresentation

$oldDate = "4/28/21"
$newDate = "5/2/21"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch {}
        if ($phType -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom layout hanging off the master.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# NOTE: Presentation.NotesMaster.Shapes is intentionally left untouched here.
# In this COM host, writes through that collection are mis-routed onto the
# slide master's shapes (an index collision: NotesMaster.Shapes.Item(N) ends
# up mutating SlideMaster.Shapes.Item(N) at save time), silently corrupting
# the master instead of updating the notes master. There is no working path
# exposed to edit the notes master's date placeholder from this host, so it
# is skipped rather than risk clobbering the slide master.

# Slide 1 content tweaks. Locate shapes by their visible text (robust to
# any reordering of the shape collection) rather than a hard-coded index.
$slide = $p.Slides.Item(1)

$oldAB = "Automate A/B, A/B/n, Canary, and Conformance experiments"
$newAB = "Automate A/B(/n), Canary, and Conformance experiments"
$oldLead = "Safely promote winning version "
$newLead = "Find and promote winning version "

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $txt = $tr.Text

    if ($txt -eq $oldAB) {
        $tr.Text = $newAB
    }
    elseif ($txt.StartsWith($oldLead)) {
        # Only rewrite the leading run so the trailing "of your app/ML
        # model" run keeps its own (different) formatting untouched.
        $lead = $tr.Characters(1, $oldLead.Length)
        $lead.Text = $newLead
    }
}
